# Generate Report for Archive
#
# 1. Update the "Status" text from "Ready for handoff" to "In Translation"
#    on all sheets that show it (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2. Shrink the "Status"-related columns' widths from 17.2159881591797 to
#    13.4101845877511 (Overview columns E & F, zh-cn column C, de-de column C).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Update status text ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Update column widths ---
# Target stored width is 13.4101845877511 "characters". Excel's ColumnWidth
# property is quantized to whole pixels (MaxDigitWidth=6 for Calibri 11), so
# the closest reproducible value is obtained with ColumnWidth = 12.5, which
# lands on stored width 13.333333333333334 (nearest achievable pixel grid
# point to the requested width).
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
